# Update gh-pages to output generated at 456a3b4
# Applies cell updates to the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value = 333
$ws1.Range("F4").Value = 270
$ws1.Range("D6").Value = "南宁国际会展中心  南宁国际会展中心"
$ws1.Range("F6").Value = 3186
$ws1.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202405/hyC2ZhnZ1715826721453.jpeg"
$ws1.Range("F7").Value = 2089
$ws1.Range("F10").Value = 1184
$ws1.Range("F12").Value = 1079
$ws1.Range("F13").Value = 89

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value = 333
$ws4.Range("F4").Value = 270
$ws4.Range("D6").Value = "南宁国际会展中心  南宁国际会展中心"
$ws4.Range("F6").Value = 3186
$ws4.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202405/hyC2ZhnZ1715826721453.jpeg"
$ws4.Range("F7").Value = 2089
$ws4.Range("F11").Value = 1184
$ws4.Range("F13").Value = 1079
$ws4.Range("F14").Value = 89
